$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 161.4224869728088
$ws.Range("C2").Value = 55.76770927508439
$ws.Range("D2").Value = 2.089466094970703
$ws.Range("E2").Value = 3.141579706245607
$ws.Range("B3").Value = 433.2628062725067
$ws.Range("C3").Value = 137.5478104021661
$ws.Range("D3").Value = 4.608819103240966
$ws.Range("E3").Value = 1.948374591872216
$ws.Range("B4").Value = 891.3809002876282
$ws.Range("C4").Value = 145.7360696628547
$ws.Range("D4").Value = 7.821775579452515
$ws.Range("E4").Value = 0.3344061402249211
$ws.Range("B5").Value = 396.2700290679932
$ws.Range("C5").Value = 5.397329041317824
$ws.Range("D5").Value = 6.274562168121338
$ws.Range("E5").Value = 0.6201628296824676
$ws.Range("B6").Value = 788.7560674190521
$ws.Range("C6").Value = 5.878420357202714
$ws.Range("D6").Value = 7.487487220764161
$ws.Range("E6").Value = 0.4733867440628987
$ws.Range("B7").Value = 1577.682591295242
$ws.Range("C7").Value = 12.35958354394973
$ws.Range("D7").Value = 7.810253477096557
$ws.Range("E7").Value = 0.5289247424152462
$ws.Range("B8").Value = 534.3931314945221
$ws.Range("C8").Value = 7.24980997792151
$ws.Range("D8").Value = 7.1660080909729
$ws.Range("E8").Value = 0.3832394366663081
$ws.Range("B9").Value = 1043.80263209343
$ws.Range("C9").Value = 6.852773957842119
$ws.Range("D9").Value = 7.343221998214721
$ws.Range("E9").Value = 0.6997019563286839
$ws.Range("B10").Value = 2071.147108125687
$ws.Range("C10").Value = 16.80539605895351
$ws.Range("D10").Value = 8.805324649810791
$ws.Range("E10").Value = 0.3610204129538772
$ws.Range("B11").Value = 384.4225076198578
$ws.Range("C11").Value = 5.232459231192113
$ws.Range("D11").Value = 7.629725456237793
$ws.Range("E11").Value = 0.6172351791407384
$ws.Range("B12").Value = 747.8065921783448
$ws.Range("C12").Value = 16.89020334314273
$ws.Range("D12").Value = 8.409038209915161
$ws.Range("E12").Value = 0.2801360180783285
$ws.Range("B13").Value = 1481.346619939804
$ws.Range("C13").Value = 26.00648650946575
$ws.Range("D13").Value = 7.798697423934937
$ws.Range("E13").Value = 0.9735022665631949
$ws.Range("B14").Value = 552.8581718921662
$ws.Range("C14").Value = 16.16228369270926
$ws.Range("D14").Value = 8.011240863800049
$ws.Range("E14").Value = 0.3569548795134374
$ws.Range("B15").Value = 1088.408327770233
$ws.Range("C15").Value = 17.43759374865452
$ws.Range("D15").Value = 6.977490282058715
$ws.Range("E15").Value = 1.352000596809704
$ws.Range("B16").Value = 2162.664258575439
$ws.Range("C16").Value = 43.81784823854721
$ws.Range("D16").Value = 8.87287130355835
$ws.Range("E16").Value = 0.4161040172149783
$ws.Range("B17").Value = 720.4558026313782
$ws.Range("C17").Value = 15.71153951083405
$ws.Range("D17").Value = 7.194146823883057
$ws.Range("E17").Value = 0.4238341641824623
$ws.Range("B18").Value = 1416.47435464859
$ws.Range("C18").Value = 17.7788358890005
$ws.Range("D18").Value = 7.840663909912109
$ws.Range("E18").Value = 0.556430291069796
$ws.Range("B19").Value = 2858.879571056366
$ws.Range("C19").Value = 38.82474564596157
$ws.Range("D19").Value = 8.501192188262939
$ws.Range("E19").Value = 0.5759237248439967
$ws.Range("B20").Value = 398.8124918460846
$ws.Range("C20").Value = 5.251522608486193
$ws.Range("D20").Value = 7.253786945343018
$ws.Range("E20").Value = 0.4704588869595814
$ws.Range("B21").Value = 792.2914300918579
$ws.Range("C21").Value = 19.7299568157689
$ws.Range("D21").Value = 7.945536661148071
$ws.Range("E21").Value = 0.5291216659281549
$ws.Range("B22").Value = 1546.652203798294
$ws.Range("C22").Value = 28.18179536927087
$ws.Range("D22").Value = 8.360516834259034
$ws.Range("E22").Value = 0.5502584810651959
$ws.Range("B23").Value = 577.9067765235901
$ws.Range("C23").Value = 12.33066342569965
$ws.Range("D23").Value = 7.272259044647217
$ws.Range("E23").Value = 0.5479504410803239
$ws.Range("B24").Value = 1138.185612773895
$ws.Range("C24").Value = 21.2134232392483
$ws.Range("D24").Value = 7.922602272033691
$ws.Range("E24").Value = 0.4284979821316703
$ws.Range("B25").Value = 2263.575600004196
$ws.Range("C25").Value = 38.76192520549186
$ws.Range("D25").Value = 8.492589855194092
$ws.Range("E25").Value = 0.4436005557169488
$ws.Range("B26").Value = 757.2915437221527
$ws.Range("C26").Value = 13.9832792791435
$ws.Range("D26").Value = 8.199797248840332
$ws.Range("E26").Value = 0.141052092428875
$ws.Range("B27").Value = 1475.197072792053
$ws.Range("C27").Value = 22.39052006420526
$ws.Range("D27").Value = 7.467642593383789
$ws.Range("E27").Value = 0.9229999410580189
$ws.Range("B28").Value = 1999.117674779892
$ws.Range("C28").Value = 296.1709115229847
$ws.Range("D28").Value = 2.397410678863525
$ws.Range("E28").Value = 1.671764377982837
